$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# "filter out 'and' and ','" -- the industry-name -> id matcher was
# re-run after stripping the literal word "and" and comma characters
# from the match text, which changes a lot of previously-defaulted
# matches (old default id 94) to their real target id, and leaves a
# few rows with no match at all (column B cleared).
# ---------------------------------------------------------------------

$updates = @{
  1 = 35
  10 = 32
  18 = 85
  19 = 35
  20 = 9
  22 = 35
  24 = 35
  26 = 35
  44 = 6
  45 = 20
  59 = 42
  62 = 79
  63 = 79
  64 = 79
  72 = 85
  82 = 35
  87 = 35
  100 = 35
  103 = 35
  122 = 89
  123 = 89
  125 = 92
  130 = 35
  133 = 90
  134 = 60
  136 = 37
  137 = 37
  143 = 138
  145 = 8
  148 = 138
  152 = 138
  154 = 90
  160 = 35
  167 = 51
  168 = 51
  187 = 22
  200 = 42
  208 = 35
  209 = 35
  213 = 69
  216 = 138
  219 = 35
  231 = 113
  241 = 70
  244 = 35
  245 = 29
  247 = 29
  254 = 41
  255 = 35
  256 = 42
  259 = 64
  271 = 43
  275 = 4
  280 = 22
  293 = 80
  304 = 84
  305 = 35
  308 = 25
  309 = 25
  311 = 90
  314 = 138
  315 = 92
  320 = 25
  321 = 90
  340 = 35
  349 = 61
  350 = 61
  351 = 42
  352 = 24
  356 = 35
  358 = 22
  376 = 42
  382 = 54
  393 = 27
  394 = 85
  396 = 42
  398 = 27
  399 = 89
  401 = 90
  402 = 73
  405 = 28
  406 = 27
  408 = 42
  409 = 27
}
foreach ($row in $updates.Keys) {
  $ws.Cells.Item($row, 2).Value = $updates[$row]
}

$clears = @(14, 23, 29, 48, 55, 79, 85, 98, 127, 129, 132, 138, 153, 184, 185, 186, 228, 229, 253, 268, 287, 300, 328, 375, 381, 387, 418)
foreach ($row in $clears) {
  $ws.Cells.Item($row, 2).ClearContents()
}

# --- Summary row: count of industries that now have a mapped id ------
$ws.Range("B419").Formula = "=COUNTA(B1:B418)"
$ws.Range("B419").Font.Bold = $true

# --- Column layout: drop column B's custom width, widen column A -----
$ws.Columns("A").ColumnWidth = 50.625

# --- Reset the view: no more scrolled-down / mid-sheet selection ------
$ws.Range("B1").Select()

